# Remove the trailing "Ver no Jupiter..." and copyright boilerplate
# paragraphs (plus the now-redundant blank paragraph that used to
# separate them from the page-break paragraph), as produced by the
# latest site build.

$d = $word.ActiveDocument

# Locate the paragraph that starts with "Ver no Jupiter ..." by
# searching for its text.
$found = $d.Content
$found.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, `
                     $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $found.Start

# Translate that position into a paragraph index so we can grab a
# stable handle on the paragraph objects involved.
$startIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $startPos) {
        $startIdx = $i
        break
    }
}

if ($startIdx -gt 0) {
    # Paragraph sequence being removed:
    #   [startIdx]     "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   [startIdx + 1] "© 2020 . Contact: ... Creative Commons Attribution"
    #   [startIdx + 2] (blank paragraph that trails the copyright line)
    $firstPara = $d.Paragraphs.Item($startIdx)
    $lastPara  = $d.Paragraphs.Item($startIdx + 2)

    $deleteRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)
    $deleteRange.Delete()
}
